$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 84: the F84 label was actually for the CRM opened 2019 batch (fix mislabeling) ---
$ws.Range("F84").Value = "New CRM opened 12/11/2019"

# --- Row 85: new rerun data (day 2 of reruns from OR tide pools) ---
$ws.Range("A85").Value = 43819
# Copy the date number-format (style) from an existing date cell so A85 matches A2..A84 (s="1")
$ws.Range("A2").Copy()
$ws.Range("A85").PasteSpecial(-4122)

$ws.Range("B85").Value = 2196.8988871280399
$ws.Range("C85").Value = 2207.0300000000002
$ws.Range("D85").Formula = "=100*(B85-C85)/C85"
$ws.Range("E85").Value = 169
$ws.Range("F85").Value = "New CRM opened 12/11/2019"

# --- Update the selected/active cell shown when the workbook was last saved ---
$ws.Range("F88").Select()
